# Kingdom Hearts Recoded workbook - "Locations" sheet update
# - Fix SoA2 display text (drop the "or Destiny Islands..." tail)
# - Fix DestinySecret display text (drop the "Destiny Islands: " prefix)
# - Remove the DISectorF2 row entirely (duplicate/obsolete entry)
# - Fix DestinyStorm display text (drop the "Destiny Islands: " prefix)
# - Append two new location rows for the Olympus Coliseum world
#   (OlympusOutside / OlympusVestibule)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Locations")

# Shorten SoA2's display string
$ws.Range("E4").Value = "Station of Awakening"

# Shorten DestinySecret's display string
$ws.Range("E6").Value = "Secret Place"

# Delete the whole "DISectorF2" row (row 7); everything below shifts up one row
$ws.Rows(7).Delete()

# After the shift, row 7 is now "DestinyStorm" -- shorten its display string too
$ws.Range("E7").Value = "Storm-tossed Island"

# Add the two new Olympus Coliseum locations in the next free rows (31 & 32,
# which already carry sequential IDs 30 / 31 in column A from the template)
$ws.Range("B31").Value = "0x1d"
$ws.Range("C31").Value = "0x3"
$ws.Range("D31").Value = "OlympusOutside"
$ws.Range("E31").Value = "Outside the Coliseum"

$ws.Range("B32").Value = "0x1e"
$ws.Range("C32").Value = "0x3"
$ws.Range("D32").Value = "OlympusVestibule"
$ws.Range("E32").Value = "Vestibule"

# Leave the selection where the author last left it while typing the new rows
$ws.Range("E36").Select()
